$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1697247706422018
$ws.Range("C2").Value = 0.6100917431192661
$ws.Range("J2").Value = 0.004587155963302753
$ws.Range("P2").Value = 0.1330275229357798
$ws.Range("S2").Value = 0.08256880733944955

# Row 3
$ws.Range("B3").Value = 0.02142857142857143
$ws.Range("C3").Value = 0.05
$ws.Range("J3").Value = 0.007142857142857143
$ws.Range("P3").Value = 0.7142857142857143
$ws.Range("S3").Value = 0.2071428571428572

# Row 4
$ws.Range("J4").Value = 0.02083333333333333
$ws.Range("P4").Value = 0.8125
$ws.Range("S4").Value = 0.1666666666666667

# Row 6
$ws.Range("B6").Value = 0.03365384615384615
$ws.Range("D6").Value = 0.01923076923076923
$ws.Range("E6").Value = 0.009615384615384616
$ws.Range("F6").Value = 0.09615384615384616
$ws.Range("J6").Value = 0.2211538461538461
$ws.Range("O6").Value = 0.01442307692307692
$ws.Range("Q6").Value = 0.1682692307692308
$ws.Range("R6").Value = 0.08653846153846154
$ws.Range("S6").Value = 0.3509615384615384

# Row 7
$ws.Range("B7").Value = 0.0776255707762557
$ws.Range("D7").Value = 0.0045662100456621
$ws.Range("F7").Value = 0.0410958904109589
$ws.Range("J7").Value = 0.136986301369863
$ws.Range("O7").Value = 0.0091324200913242
$ws.Range("Q7").Value = 0.1963470319634703
$ws.Range("R7").Value = 0.0730593607305936
$ws.Range("S7").Value = 0.4611872146118721

# Row 8
$ws.Range("B8").Value = 0.07380952380952381
$ws.Range("D8").Value = 0.01666666666666667
$ws.Range("F8").Value = 0.04523809523809524
$ws.Range("J8").Value = 0.1238095238095238
$ws.Range("O8").Value = 0.01904761904761905
$ws.Range("Q8").Value = 0.1642857142857143
$ws.Range("R8").Value = 0.1119047619047619
$ws.Range("S8").Value = 0.4452380952380952

# Row 9
$ws.Range("B9").Value = 0.06779661016949153
$ws.Range("D9").Value = 0.01271186440677966
$ws.Range("F9").Value = 0.05508474576271186
$ws.Range("J9").Value = 0.1398305084745763
$ws.Range("O9").Value = 0.02542372881355932
$ws.Range("Q9").Value = 0.211864406779661
$ws.Range("R9").Value = 0.09322033898305085
$ws.Range("S9").Value = 0.3940677966101695

# Row 10
$ws.Range("B10").Value = 0.08122205663189269
$ws.Range("D10").Value = 0.02757078986587183
$ws.Range("F10").Value = 0.07004470938897168
$ws.Range("J10").Value = 0.1363636363636364
$ws.Range("O10").Value = 0.009687034277198211
$ws.Range("Q10").Value = 0.2116244411326378
$ws.Range("R10").Value = 0.09910581222056632
$ws.Range("S10").Value = 0.364381520119225

# Row 11
$ws.Range("G11").Value = 0.1725239616613418
$ws.Range("J11").Value = 0.03833865814696485
$ws.Range("K11").Value = 0.2044728434504792
$ws.Range("L11").Value = 0.5654952076677316
$ws.Range("S11").Value = 0.01916932907348243

# Row 12
$ws.Range("G12").Value = 0.7365591397849462
$ws.Range("J12").Value = 0.1935483870967742
$ws.Range("K12").Value = 0.005376344086021506
$ws.Range("L12").Value = 0.02688172043010753
$ws.Range("S12").Value = 0.03763440860215054

# Row 13
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2888888888888889
$ws.Range("S13").Value = 0.04444444444444445

# Row 14
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333

# Row 15
$ws.Range("F15").Value = 0.01142857142857143
$ws.Range("H15").Value = 0.16
$ws.Range("I15").Value = 0.09142857142857143
$ws.Range("J15").Value = 0.3657142857142857
$ws.Range("K15").Value = 0.04571428571428571
$ws.Range("M15").Value = 0.01142857142857143
$ws.Range("O15").Value = 0.04
$ws.Range("S15").Value = 0.2742857142857143

# Row 16
$ws.Range("F16").Value = 0.01829268292682927
$ws.Range("H16").Value = 0.1890243902439024
$ws.Range("I16").Value = 0.1219512195121951
$ws.Range("J16").Value = 0.4085365853658536
$ws.Range("K16").Value = 0.1219512195121951
$ws.Range("M16").Value = 0.03048780487804878
$ws.Range("O16").Value = 0.02439024390243903
$ws.Range("S16").Value = 0.08536585365853659

# Row 17
$ws.Range("F17").Value = 0.00631578947368421
$ws.Range("H17").Value = 0.16
$ws.Range("I17").Value = 0.12
$ws.Range("J17").Value = 0.4252631578947368
$ws.Range("K17").Value = 0.1157894736842105
$ws.Range("M17").Value = 0.01263157894736842
$ws.Range("O17").Value = 0.04842105263157895
$ws.Range("S17").Value = 0.1115789473684211

# Row 18
$ws.Range("F18").Value = 0.01265822784810127
$ws.Range("H18").Value = 0.2025316455696203
$ws.Range("I18").Value = 0.1054852320675106
$ws.Range("J18").Value = 0.4261603375527426
$ws.Range("K18").Value = 0.1012658227848101
$ws.Range("M18").Value = 0.008438818565400843
$ws.Range("N18").Value = 0.004219409282700422
$ws.Range("O18").Value = 0.04219409282700422
$ws.Range("S18").Value = 0.0970464135021097

# Row 19
$ws.Range("F19").Value = 0.0168
$ws.Range("H19").Value = 0.1896
$ws.Range("I19").Value = 0.0936
$ws.Range("J19").Value = 0.4112
$ws.Range("K19").Value = 0.1112
$ws.Range("M19").Value = 0.024
$ws.Range("N19").Value = 0.0024
$ws.Range("O19").Value = 0.056
$ws.Range("S19").Value = 0.09520000000000001
